# Add a new "2021" column (P) to the 4.2.2 participation-rate table,
# mirroring the formatting already used for the "2020" column (O).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell (year) - copy O4's format onto P4, then set the new year.
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P4").Value = 2021

# Data cell (rate) - copy O5's format onto P5, then set the new value.
$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P5").Value = 80.900000000000006

$excel.CutCopyMode = $false

# Match the author's final selection.
$ws.Range("N10").Select()
